$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The opening paragraph currently holds the "_GoBack" bookmark followed by
# the "Dear Data Science Team Leader," run, both inside the same
# paragraph. Split that into two paragraphs: the salutation text moves to
# its own (first) paragraph, and the (now empty) bookmark paragraph
# follows it.

$p1 = $d.Paragraphs.Item(1)
$rng1 = $p1.Range

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Dear Data Science Team Leader,</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$rng1.InsertXML($xml1)

# --- Change 2 -----------------------------------------------------------
# Remove the "As a reminder, ..." reminder sentence entirely (its run is
# deleted, leaving the empty paragraph with its original formatting in
# place).

$reminder = "As a reminder, the client indicated that they wanted to know the following: " +
    [char]0x201C + "How to better stock the items that they sell." + [char]0x201D + " "

$findRng = $d.Content
$found = $findRng.Find.Execute($reminder, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    [void]$findRng.Delete()
}
